# Update min_price (column D) and park-auto-sm.ru_price (column F)
# for the rows whose prices changed in this data refresh.
# Each entry: RowNumber = NewPrice
$updates = @{
    9  = 972000
    12 = 1001940
    16 = 1956000
    20 = 2495940
    21 = 1307940
    23 = 1421940
    24 = 1967940
    25 = 2777940
    26 = 1319000
    27 = 1631940
    40 = 2490000
    67 = 1499994
    68 = 1033300
    70 = 1133994
    72 = 2303994
    73 = 1961994
    79 = 1787000
    80 = 2069400
    82 = 1709400
    88 = 1289400
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in $updates.Keys) {
    $price = $updates[$row]
    $ws.Range("D$row").Value = $price
    $ws.Range("F$row").Value = $price
}
